$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format column D (Price) as Text first so numeric-looking values like
# "219.11" or "0.523" are written as literal text, matching the source
# workbook which stores every Price/Volume cell as an inline string.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.094.44"
$ws.Range("E2").Value = "  +3.71%  "
$ws.Range("D3").Value = "1.728.29"
$ws.Range("E3").Value = "  +2.85%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "219.11"
$ws.Range("E5").Value = "  +1.81%  "
$ws.Range("D6").Value = "0.523"
$ws.Range("E6").Value = "  +1.22%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "24.15"
$ws.Range("E8").Value = "  +13.38%  "
$ws.Range("D9").Value = "0.265"
$ws.Range("E9").Value = "  +3.34%  "
$ws.Range("E10").Value = "  +1.71%  "
$ws.Range("D11").Value = "0.0899"
$ws.Range("E11").Value = "  +1.68%  "
$ws.Range("D12").Value = "1.971.86"
$ws.Range("E12").Value = "  +2.89%  "
$ws.Range("D13").Value = "1.727.00"
$ws.Range("E13").Value = "  +2.90%  "
$ws.Range("E14").Value = "  +3.70%  "
$ws.Range("E15").Value = "  +5.02%  "
$ws.Range("D16").Value = "67.56"
$ws.Range("E16").Value = "  +2.12%  "
$ws.Range("D17").Value = "28.055.95"
$ws.Range("E17").Value = "  +3.58%  "
$ws.Range("D18").Value = "243.21"
$ws.Range("E18").Value = "  +2.44%  "
$ws.Range("E19").Value = "  +1.81%  "
$ws.Range("D20").Value = "7.92"
$ws.Range("E20").Value = "  -2.86%  "
$ws.Range("E22").Value = "  +3.64%  "
$ws.Range("E23").Value = "  +3.93%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").Value = "148.92"
$ws.Range("E25").Value = "  +1.50%  "
$ws.Range("D26").Value = "7.53"
$ws.Range("E26").Value = "  +4.10%  "
$ws.Range("D27").Value = "16.75"
$ws.Range("E27").Value = "  +2.65%  "
$ws.Range("E28").Value = "  +1.46%  "
$ws.Range("D30").Value = "0.0510"
$ws.Range("E30").Value = "  +2.33%  "
$ws.Range("E31").Value = "  +2.10%  "
$ws.Range("E32").Value = "  +2.46%  "
$ws.Range("D33").Value = "1.496.17"
$ws.Range("E33").Value = "  -3.81%  "
$ws.Range("D34").Value = "3.28"
$ws.Range("E34").Value = "  +2.63%  "
$ws.Range("D35").Value = "1.67"
$ws.Range("E35").Value = "  -1.88%  "
$ws.Range("D36").Value = "0.610"
$ws.Range("E36").Value = "  +1.26%  "
$ws.Range("D37").Value = "0.956"
$ws.Range("E37").Value = "  +1.81%  "
$ws.Range("E38").Value = "  +1.13%  "
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("E40").Value = "  +0.85%  "
$ws.Range("D41").Value = "70.75"
$ws.Range("E41").Value = "  +3.67%  "
$ws.Range("E42").Value = "  +3.75%  "
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("E44").Value = "  +2.29%  "
$ws.Range("D45").Value = "1.876.27"
$ws.Range("E45").Value = "  +2.82%  "
$ws.Range("D46").Value = "0.800"
$ws.Range("E46").Value = "  +1.99%  "
$ws.Range("D47").Value = "1.79"
$ws.Range("E47").Value = "  +13.89%  "
$ws.Range("D48").Value = "91.19"
$ws.Range("E48").Value = "  +0.40%  "
$ws.Range("E49").Value = "  +4.09%  "
$ws.Range("D50").Value = "8.26"
$ws.Range("E50").Value = "  +2.96%  "
$ws.Range("E51").Value = "  +0.66%  "

# Restore the default "Normal" style on column D so no stray cell-level
# style index lingers on cells that did not have one before.
$priceRange.Style = "Normal"
